# Fix the "cold space SOS" scene-name string: replace the Unicode en-dash
# (U+2013) with a plain ASCII hyphen-minus, as called out in the commit
# message ("update template to fix scene name with unicode").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "cold space SOS - Scene01"

# Move the active selection to K20 (matches the sheetView/selection saved
# in the workbook at commit time).
[void]$ws.Range("K20").Select()

# The column widths were re-measured (tiny sub-pixel shrink across every
# custom column) when the workbook was last saved. Reproduce the new
# widths for the explicitly-sized columns A:L; leave the trailing
# default-width column block (M onward) alone so the sheet keeps a single
# merged <col> run there instead of splitting into per-column entries.
$ws.Columns.Item(1).ColumnWidth = 11.333333333333332
$ws.Columns.Item(2).ColumnWidth = 8.5
$ws.Columns.Item(3).ColumnWidth = 22.833333333333336
$ws.Columns.Item(4).ColumnWidth = 16.333333333333336
$ws.Columns.Item(5).ColumnWidth = 14.666666666666666
$ws.Columns.Item(6).ColumnWidth = 14.333333333333332
$ws.Columns.Item(7).ColumnWidth = 17.333333333333336
$ws.Columns.Item(8).ColumnWidth = 15.0
$ws.Columns.Item(9).ColumnWidth = 14.333333333333332
$ws.Columns.Item(10).ColumnWidth = 27.5
$ws.Columns.Item(11).ColumnWidth = 20.666666666666668
$ws.Columns.Item(12).ColumnWidth = 14.333333333333332
